# Add a "Seller Slug" column before the existing "Review Count" column
# (old D/E/F -> new E/F/G), then fill it in with the seller slug parsed
# out of each row's NOTHS URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; this shifts Review Count/NOTHS URL/Feefo URL
# one column to the right (D->E, E->F, F->G) and carries over the header style.
$ws.Range("D1").EntireColumn.Insert()

# New header
$ws.Range("D1").Value = "Seller Slug"

# Seller slug values, derived from the seller segment of each row's NOTHS URL
# (now in column F after the insert), row by row for rows 2..22.
$slugs = @{
    2  = "ellieellie"
    3  = "thechucklingcheesecompany"
    4  = "dibor"
    5  = "gaamaa"
    6  = "oakdenedesigns"
    7  = "myposhshop"
    8  = "lovetreedesign"
    9  = "madewithlovecardboutique"
    10 = "songsofinkandsteel"
    11 = "thegourmetchocolatepizzaco"
    12 = "qwertybeerbox"
    13 = "thealphabetgiftshop"
    14 = "ladedaliving"
    15 = "hurleyburleyman"
    16 = "joybycorrinesmith"
    17 = "therusticdish"
    18 = "lisaangeljewellery"
    19 = "theforestandco"
    20 = "hurleyburleyman"
    21 = "alphabetinteriors"
    22 = "dibor"
}

foreach ($row in $slugs.Keys) {
    $ws.Cells.Item($row, 4).Value = $slugs[$row]
}
